# Generate Report for handoff
#
# The "bb511d85-..." localization unit was re-handed-off (new handoff
# timestamp), which flips its status from "Handed back: in sync with en-US"
# to "Ready for handoff" and re-sorts it after the "c79106f3-..." unit on
# every report sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$READY = "Ready for handoff"
$HANDED_BACK = "Handed back: in sync with en-US"

$BB = "bb511d85-f275-4989-80c1-b976714011e0"
$CF = "c79106f3-e814-4dc5-92c8-b94d92bfafdb"

# ---------------------------------------------------------------------
# Helper: set a cell's text value and, if that cell carries a hyperlink,
# update the hyperlink's visible display text to match (done via foreach
# so the existing <hyperlink> element is edited in place instead of a
# duplicate being appended).
# ---------------------------------------------------------------------
function Set-CellAndLink($ws, $addr, $text) {
    $ws.Range($addr).Value2 = $text
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

# ========================= Overview sheet =============================
$ov = $wb.Worksheets.Item("Overview")

Set-CellAndLink $ov '$A$2' ($CF + ".md")
Set-CellAndLink $ov '$A$3' ($BB + ".md")

$ov.Range("B3").Value2 = $READY
$ov.Range("C3").Value2 = $READY

# ====================== zh-cn / de-de detail sheets =====================
$langs = @(
    @{ Name = "zh-cn"; Ext = "zh-cn.xlf"; BbHandoff = "2016-01-17 10:11:03"; BbNewHandoff = "2016-01-17 10:12:49"; CfHandoff = "2016-01-17 10:11:03"; CfHandback = "2016-01-17 10:11:47"; BbHandback = "2016-01-17 10:11:47" },
    @{ Name = "de-de"; Ext = "de-de.xlf"; BbHandoff = "2016-01-17 10:11:15"; BbNewHandoff = "2016-01-17 10:13:00"; CfHandoff = "2016-01-17 10:11:15"; CfHandback = "2016-01-17 10:12:06"; BbHandback = "2016-01-17 10:12:06" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Name)

    $cfHandoffFile = $CF + ".4bf28bf780fa09b40203412d942fdfdc6bd5b109." + $lang.Ext
    $bbHandoffFile = $BB + ".20500ee177d492de74dded2302cf4d0661315451." + $lang.Ext

    # Row 2 now holds the c79106f3 unit (unchanged content, just moved up)
    Set-CellAndLink $ws '$A$2' ($CF + ".md")
    $ws.Range("B2").Value2 = $HANDED_BACK
    Set-CellAndLink $ws '$C$2' $cfHandoffFile
    $ws.Range("D2").Value2 = $lang.CfHandoff
    Set-CellAndLink $ws '$E$2' ($CF + ".md")
    Set-CellAndLink $ws '$F$2' $cfHandoffFile
    $ws.Range("G2").Value2 = $lang.CfHandback
    $ws.Range("H2").Value2 = "Include"

    # Row 3 now holds the bb511d85 unit, with an updated handoff datetime
    # and a status bumped to "Ready for handoff"
    Set-CellAndLink $ws '$A$3' ($BB + ".md")
    $ws.Range("B3").Value2 = $READY
    Set-CellAndLink $ws '$C$3' $bbHandoffFile
    $ws.Range("D3").Value2 = $lang.BbNewHandoff
    Set-CellAndLink $ws '$E$3' ($BB + ".md")
    Set-CellAndLink $ws '$F$3' $bbHandoffFile
    $ws.Range("G3").Value2 = $lang.BbHandback
    $ws.Range("H3").Value2 = "Include"
}
